$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above row 2 -- this shifts the old rows 2-7 down to 3-8
#    and keeps all their original formatting/content intact.
$ws.Rows(2).Insert()

# 2. Fill in the new subtitle row (row 2) with the "(in percent)" translations.
$ws.Range("A2").Value = "(пайыз менен)"
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"

# 3. Style the new subtitle cells: italic, 8pt Times New Roman, centered.
$sub = $ws.Range("A2:C2")
$sub.Font.Name = "Times New Roman"
$sub.Font.Size = 8
$sub.Font.Italic = $true
$sub.HorizontalAlignment = -4108
$sub.VerticalAlignment = -4108

# 4. Add the new 2023 column (column O) mirroring the 2022 column's data/style.
$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 6.18

# 5. Row-height tweaks that came along with the edit.
$ws.Rows(1).RowHeight = 60.75
$ws.Rows(5).RowHeight = 53.25

# 6. Column width tweak: columns A-C now share one uniform width.
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 39.307291666666664

# 7. Leave the cursor on A1 (matches a "clean" selection state).
$ws.Range("A1").Select()
